$d = $word.ActiveDocument

# Replace the identifier text (merges the two runs - text + trailing space - into one)
$d.Content.Find.Execute("**ID__AFFARS_pgi_5342_topic_5__ID** ", $true, $false, $false, $false, $false, $true, 1, $false, "**ID__AFFARS_AFMC_PGI_5342_1503_90__ID**", 2)

# Update the first paragraph's formatting: add a paragraph border (space-only) and change left indent
$p = $d.Paragraphs(1)
$p.Format.Borders.DistanceFromTop = 5
$p.Format.Borders.DistanceFromLeft = 5
$p.Format.Borders.DistanceFromBottom = 5
$p.Format.Borders.DistanceFromRight = 5
$p.Format.LeftIndent = 11.25
